# Applies the profit-recalculation data refresh (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, per the scheduled
# runner's commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 348.34784
$ws.Range("I28").Value = 334.1
$ws.Range("J28").Value = 443.33334
$ws.Range("K28").Value = 334.1
$ws.Range("L28").Value = 443.33334
$ws.Range("M28").Value = 150.9
$ws.Range("N28").Value = -1413.33334

$ws.Range("H98").Value = 7867.1665
$ws.Range("I98").Value = 11051
$ws.Range("J98").Value = 1499.5
$ws.Range("K98").Value = 11051
$ws.Range("L98").Value = 1499.5
$ws.Range("M98").Value = -9553
$ws.Range("N98").Value = -4495.5

$ws.Range("H116").Value = 3315.56
$ws.Range("I116").Value = 3130.8333
$ws.Range("J116").Value = 3486.077
$ws.Range("K116").Value = 3130.8333
$ws.Range("L116").Value = 3486.077
$ws.Range("M116").Value = 311.1667000000002
$ws.Range("N116").Value = -10370.077

$ws.Range("H122").Value = 7867.1665
$ws.Range("I122").Value = 11051
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 33153
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -30703
$ws.Range("N122").Value = -9398.5

$ws.Range("H135").Value = 1160.125
$ws.Range("I135").Value = 524
$ws.Range("J135").Value = 3916.6667
$ws.Range("K135").Value = 4716
$ws.Range("L135").Value = 35250.0003
$ws.Range("M135").Value = -2181
$ws.Range("N135").Value = -40320.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5001196.5
$ws.Range("I2").Value = 8621563
$ws.Range("J2").Value = 1642.8096
$ws.Range("K2").Value = 8621563
$ws.Range("L2").Value = 1642.8096
$ws.Range("M2").Value = -8621450
$ws.Range("N2").Value = -1868.8096

$ws.Range("H4").Value = 1007
$ws.Range("I4").Value = 509.8
$ws.Range("J4").Value = 2250
$ws.Range("K4").Value = 509.8
$ws.Range("L4").Value = 2250
$ws.Range("M4").Value = -393.8
$ws.Range("N4").Value = -2482

$ws.Range("H5").Value = 599.5
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 999
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -1223

$ws.Range("H19").Value = 24903.4
$ws.Range("I19").Value = 1499.6666
$ws.Range("J19").Value = 60009
$ws.Range("K19").Value = 1499.6666
$ws.Range("L19").Value = 60009
$ws.Range("M19").Value = -1270.6666
$ws.Range("N19").Value = -60467

$ws.Range("H25").Value = 33872.668
$ws.Range("I25").Value = 1960
$ws.Range("J25").Value = 73763.5
$ws.Range("K25").Value = 1960
$ws.Range("L25").Value = 73763.5
$ws.Range("M25").Value = -1558
$ws.Range("N25").Value = -74567.5

$ws.Range("H32").Value = 1710.15
$ws.Range("I32").Value = 1710.15
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1710.15
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1423.15
$ws.Range("N32").ClearContents()

$ws.Range("H113").Value = 34000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 34000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 34000
$ws.Range("N113").Value = -42678

$ws.Range("H116").Value = 5001196.5
$ws.Range("I116").Value = 8621563
$ws.Range("J116").Value = 1642.8096
$ws.Range("K116").Value = 8621563
$ws.Range("L116").Value = 1642.8096
$ws.Range("M116").Value = -8619269
$ws.Range("N116").Value = -6230.809600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5001196.5
$ws.Range("I3").Value = 8621563
$ws.Range("J3").Value = 1642.8096
$ws.Range("K3").Value = 8621563
$ws.Range("L3").Value = 1642.8096
$ws.Range("M3").Value = -8621449
$ws.Range("N3").Value = -1870.8096

$ws.Range("H4").Value = 599.5
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 999
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 999
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -1229

$ws.Range("H15").Value = 83340.336
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 83340.336
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 83340.336
$ws.Range("N15").Value = -83794.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 79.71429000000001
$ws.Range("I7").Value = 79.8
$ws.Range("J7").Value = 79.666664
$ws.Range("K7").Value = 79.8
$ws.Range("L7").Value = 79.666664
$ws.Range("M7").Value = 33.2
$ws.Range("N7").Value = -305.666664

$ws.Range("H22").Value = 341.55554
$ws.Range("I22").Value = 296.42856
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 296.42856
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = 53.57144
$ws.Range("N22").Value = -1199.5

$ws.Range("H52").Value = 51250
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 51250
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 51250
$ws.Range("N52").Value = -51838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 17006.834
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 20398.2
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 61194.60000000001
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -61474.60000000001

$ws.Range("H47").Value = 1457.1
$ws.Range("I47").Value = 143.83333
$ws.Range("J47").Value = 3427
$ws.Range("K47").Value = 431.49999
$ws.Range("L47").Value = 10281
$ws.Range("M47").Value = -0.4999899999999684
$ws.Range("N47").Value = -11143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 90.28125
$ws.Range("I2").Value = 98.478264
$ws.Range("J2").Value = 69.333336
$ws.Range("K2").Value = 98.478264
$ws.Range("L2").Value = 69.333336
$ws.Range("M2").Value = 14.521736
$ws.Range("N2").Value = -295.333336

$ws.Range("H47").Value = 70031
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 70031
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 70031
$ws.Range("N47").Value = -71167

$ws.Range("H52").Value = 42516.5
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 42516.5
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 42516.5
$ws.Range("N52").Value = -43034.5

$ws.Range("H102").Value = 43405.08
$ws.Range("I102").Value = 2602.182
$ws.Range("J102").Value = 75464.5
$ws.Range("K102").Value = 2602.182
$ws.Range("L102").Value = 75464.5
$ws.Range("M102").Value = -980.1819999999998
$ws.Range("N102").Value = -78708.5

$ws.Range("H122").Value = 3427.7036
$ws.Range("I122").Value = 2784.3125
$ws.Range("J122").Value = 4363.5454
$ws.Range("K122").Value = 8352.9375
$ws.Range("L122").Value = 13090.6362
$ws.Range("M122").Value = -5902.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 55556492
$ws.Range("I16").Value = 76924136
$ws.Range("J16").Value = 625.4
$ws.Range("K16").Value = 76924136
$ws.Range("L16").Value = 625.4
$ws.Range("M16").Value = -76923966
$ws.Range("N16").Value = -965.4

$ws.Range("H53").Value = 5180
$ws.Range("I53").Value = 3950
$ws.Range("J53").Value = 6000
$ws.Range("K53").Value = 3950
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = -3432
$ws.Range("N53").Value = -7036

$ws.Range("H54").Value = 14000.667
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14000.667
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14000.667
$ws.Range("N54").Value = -15288.667

$ws.Range("H132").Value = 7663.385
$ws.Range("I132").Value = 3093.5
$ws.Range("J132").Value = 9694.444
$ws.Range("K132").Value = 9280.5
$ws.Range("L132").Value = 29083.332
$ws.Range("M132").Value = -6750.5
$ws.Range("N132").Value = -34143.33199999999

$ws.Range("H136").Value = 2463.25
$ws.Range("I136").Value = 2223.5483
$ws.Range("J136").Value = 3288.889
$ws.Range("K136").Value = 6670.644899999999
$ws.Range("L136").Value = 9866.667000000001
$ws.Range("M136").Value = -4120.644899999999
$ws.Range("N136").Value = -14966.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18000
$ws.Range("N51").Value = -19020

$ws.Range("H53").Value = 19000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 19000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 19000
$ws.Range("N53").Value = -20214

$ws.Range("H81").Value = 1830.6364
$ws.Range("I81").Value = 1833.6666
$ws.Range("J81").Value = 1827
$ws.Range("K81").Value = 3667.3332
$ws.Range("L81").Value = 3654
$ws.Range("M81").Value = -2606.3332
$ws.Range("N81").Value = -5776

$ws.Range("H84").Value = 1830.6364
$ws.Range("I84").Value = 1833.6666
$ws.Range("J84").Value = 1827
$ws.Range("K84").Value = 18336.666
$ws.Range("L84").Value = 18270
$ws.Range("M84").Value = -13032.666
$ws.Range("N84").Value = -28878

$ws.Range("H107").Value = 367
$ws.Range("I107").Value = 380.4
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1141.2
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 778.8000000000002
$ws.Range("N107").Value = -4740

$ws.Range("H113").Value = 882.7895
$ws.Range("I113").Value = 414.3846
$ws.Range("J113").Value = 1897.6666
$ws.Range("K113").Value = 1243.1538
$ws.Range("L113").Value = 5692.9998
$ws.Range("M113").Value = 926.8462
$ws.Range("N113").Value = -10032.9998

$ws.Range("H122").Value = 2807.1667
$ws.Range("I122").Value = 1743
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 5229
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -2779
